$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three test e-mail addresses used double underscores; the updated
# test data replaces them with double dots. The same shared string backs
# both the "Email" (column A) and "Confirm Email" (column B) cells, so
# updating both keeps them pointing at one shared string entry.
$ws.Range("A2").Value = "manvir..singh1@gmail.com"
$ws.Range("B2").Value = "manvir..singh1@gmail.com"

$ws.Range("A3").Value = "manvir..singh12@gmail.com"
$ws.Range("B3").Value = "manvir..singh12@gmail.com"

$ws.Range("A4").Value = "manvir..singh123@gmail.com"
$ws.Range("B4").Value = "manvir..singh123@gmail.com"

# Row 2 (A2/B2) previously used a hyperlink font without the underline
# (distinct from the rest of the hyperlinked cells). Normalize it to the
# same underlined hyperlink look used by rows 3 & 4.
$ws.Range("A2:B2").Font.Underline = 2

# Leave the active selection on B2:B4, matching the final state of the
# edit session.
$ws.Range("B2:B4").Select() | Out-Null
